$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Default")

$ws.Range("A6").Value = "T"
$ws.Range("B6").Value = "KeyManager"

$ws.Range("B6").Select()
